$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$wins = 67
$losses = 94
$ties = 0

for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}
